# Calcolo Utilizzazione.xlsx - "Correzioni e Simulation Time"
# Extends the simulation table from 11 to 26 rows, updates the input
# parameters in row 4 (G4/H4/I4/K4/L4 values, J4 becomes a formula =G4),
# resizes the "Tabella1" table to match, and moves the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new input values -------------------------------------------
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 75
$ws.Range("J4").Formula = "=G4"
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 1000

# --- Extend the simulation rows from 12 through 26 ----------------------
$ws.Range("B12:B26").Formula = "=B11+1"
$ws.Range("C12:C26").Formula = "=POWER(C`$4,B12)"
$ws.Range("D12:D26").Formula = "=C12+D11*D`$4"
$ws.Range("E12:E26").Formula = "=D12+E11*E`$4"
$ws.Range("G12:I26").Formula = "=C`$4*(`$E11/`$E12)"

# --- Resize the table to cover the new rows ------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B3:E26"))

# --- Update the view: scroll / selection ---------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$excel.Goto($ws.Range("J4:L4"))
